# Applies the "Future arm movement updates" edit to the Lessons Learned doc.
# Uses Find/Replace (wdFindContinue=1, wdReplaceOne=2) against the whole
# document content range for each distinct textual change in the diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: find failed for: $old"
    }
}

# 1) "Got IMUs now ..." -> "Have ability to attach IMUs to ..." (head IMU attachment sentence)
Replace-Text `
    "Got IMUs now so that can attach them to my head and mouth so that prop mimics my movements.  Head IMU can be attached with its rubber band and duct tape to the top of a baseball cap (see orientation markings on the IMU)" `
    "Have ability to attach IMUs to head and mouth so that prop mimics my movements.  Head IMU attached with duct tape to the top of a baseball cap (see orientation markings on the IMU)"

# 2) "Look at adjusting ..." -> "Adjusted the ..." (mouth scaling sentence, now past tense)
Replace-Text `
    "Look at adjusting the software scaling for the mouth movements to optimize it so mouth opens wider than normal mouth when normal mouth fully open" `
    "Adjusted the software scaling for the mouth movements to optimize it so prop mouth opens wider than normal mouth does"

# 3) Append new sentence about orientation setting after "... like 255 or greater."
Replace-Text `
    "the IMU heading reading suddenly jumps to large values like 255 or greater." `
    "the IMU heading reading suddenly jumps to large values like 255 or greater.  There may be an orientation setting that could solve this but simply having the IMU upside pointing up fixed the problem."

# 4) New leading space before "Implemented way to feed mic audio" paragraph
Replace-Text `
    "Implemented way to feed mic audio" `
    " Implemented way to feed mic audio"

# 5) Arm IMU count bumped from 2 to 3, plus a new clause describing what the IMUs allow
Replace-Text `
    ".  Would need another 2 IMUs for choreographing for this  (already have one but having a spare is always a good idea). It appears that should be able to have another 9 servos" `
    ".  Would need another 3 IMUs for choreographing for this  (already have one but having a spare is always a good idea) to allow arm shoulder up/down, arm shoulder side to side, and forearm up/down (a fourth for hand up/down or forearm rotate around its axis)It appears that should be able to have another 9 servos"

Write-Output "Done."
